# Grump metadata: add three new variable-metadata rows (Longhurst_Long,
# Longhurst_Short, Season) to the "vars_meta_data" sheet, directly below
# the existing rows (which end at row 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")
$ws.Activate()

# Each entry: short name / long description / sensor / spatial res /
# temporal res / discipline / visualize flag, matching columns A-H of
# the existing rows in this sheet.
$rows = @(
    @{ Row = 41; A = "Longhurst_Long";  B = "Longhurst province sample was taken in." },
    @{ Row = 42; A = "Longhurst_Short"; B = "Longhurst province sample was taken in, shortened code." },
    @{ Row = 43; A = "Season";          B = "Season sample was taken in." }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = "NA"
    $ws.Range("D$n").Value = "NA"
    $ws.Range("E$n").Value = "Irregular"
    $ws.Range("F$n").Value = "Irregular"
    $ws.Range("G$n").Value = "Biology"
    $ws.Range("H$n").Value = 1
}

$ws.Range("A41:H43").Select()
